# Applies the commit "Atualizado por script em 23-11-2023 14:45"
#  - Row 115 and Row 116 swap their match data (F:V)
#  - Row 126 and Row 130 swap their match data (F:V)
#  - A brand-new row 131 (Metalac vs Macva) is appended
#  - The sheet's dimension grows from A1:V130 to A1:V131

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 115 <-> 116 : swap "Graficar Beograd - Metalac" and
#                    "OFK Beograd - Sloboda" match rows (columns F:V)
# ---------------------------------------------------------------------
$ws.Range("F115").Value = "OFK Beograd"
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = "Sloboda"
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 1.26
$ws.Range("K115").Value = "12/11/2023 02:12"
$ws.Range("L115").Value = 1.32
$ws.Range("M115").Value = "12/11/2023 12:21"
$ws.Range("N115").Value = 4.79
$ws.Range("O115").Value = "12/11/2023 02:12"
$ws.Range("P115").Value = 4.55
$ws.Range("Q115").Value = "12/11/2023 12:21"
$ws.Range("R115").Value = 8.58
$ws.Range("S115").Value = "12/11/2023 02:12"
$ws.Range("T115").Value = 8.029999999999999
$ws.Range("U115").Value = "12/11/2023 12:21"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-sloboda/E5mu4PSG/"

$ws.Range("F116").Value = "Graficar Beograd"
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = "Metalac"
$ws.Range("I116").Value = 1
$ws.Range("J116").Value = 1.75
$ws.Range("K116").Value = "12/11/2023 02:12"
$ws.Range("L116").Value = 1.78
$ws.Range("M116").Value = "12/11/2023 12:48"
$ws.Range("N116").Value = 3.23
$ws.Range("O116").Value = "12/11/2023 02:12"
$ws.Range("P116").Value = 3.26
$ws.Range("Q116").Value = "12/11/2023 12:48"
$ws.Range("R116").Value = 4.14
$ws.Range("S116").Value = "12/11/2023 02:12"
$ws.Range("T116").Value = 4.21
$ws.Range("U116").Value = "12/11/2023 12:48"
$ws.Range("V116").Value = "https://www.betexplorer.com/football/serbia/prva-liga/graficar-beograd-metalac/f319nccq/"

# ---------------------------------------------------------------------
# Rows 126 <-> 130 : swap "RFK Novi Sad - Radnicki S. Mitrovica" and
#                    "Sloboda - FK Indjija" match rows (columns F:V)
# ---------------------------------------------------------------------
$ws.Range("F126").Value = "Sloboda"
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = "FK Indjija"
$ws.Range("I126").Value = 1
$ws.Range("J126").Value = 2.53
$ws.Range("K126").Value = "25/09/2023 07:12"
$ws.Range("L126").Value = 3.01
$ws.Range("M126").Value = "22/11/2023 12:42"
$ws.Range("N126").Value = 2.68
$ws.Range("O126").Value = "25/09/2023 07:12"
$ws.Range("P126").Value = 2.66
$ws.Range("Q126").Value = "22/11/2023 12:42"
$ws.Range("R126").Value = 2.65
$ws.Range("S126").Value = "25/09/2023 07:12"
$ws.Range("T126").Value = 2.33
$ws.Range("U126").Value = "22/11/2023 12:42"
$ws.Range("V126").Value = "https://www.betexplorer.com/football/serbia/prva-liga/sloboda-indjija/Ot0qLN2F/"

$ws.Range("F130").Value = "RFK Novi Sad"
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = "Radnicki S. Mitrovica"
$ws.Range("I130").Value = 2
$ws.Range("J130").Value = 3.46
$ws.Range("K130").Value = "26/09/2023 03:12"
$ws.Range("L130").Value = 3.16
$ws.Range("M130").Value = "22/11/2023 12:56"
$ws.Range("N130").Value = 2.98
$ws.Range("O130").Value = "26/09/2023 03:12"
$ws.Range("P130").Value = 2.86
$ws.Range("Q130").Value = "22/11/2023 12:15"
$ws.Range("R130").Value = 1.94
$ws.Range("S130").Value = "26/09/2023 03:12"
$ws.Range("T130").Value = 2.28
$ws.Range("U130").Value = "22/11/2023 12:56"
$ws.Range("V130").Value = "https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-radnicki-s-mitrovica/YFZSCaAe/"

# ---------------------------------------------------------------------
# New row 131 : Metalac vs Macva (append to end of table)
# First clone the number-formats/styles from the row above (A131 keeps
# the bold/bordered index style, E131 keeps the date/time style) then
# fill in the values.
# ---------------------------------------------------------------------
$ws.Range("A130").Copy($ws.Range("A131"))
$ws.Range("E130").Copy($ws.Range("E131"))

$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "serbia"
$ws.Range("C131").Value = "prva-liga"
$ws.Range("D131").Value = "2023-2024"
$ws.Range("E131").Value = 45253.625
$ws.Range("F131").Value = "Metalac"
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = "Macva"
$ws.Range("I131").Value = 2
$ws.Range("J131").Value = 1.99
$ws.Range("K131").Value = "26/09/2023 05:12"
$ws.Range("L131").Value = 2.11
$ws.Range("M131").Value = "23/11/2023 14:13"
$ws.Range("N131").Value = 2.78
$ws.Range("O131").Value = "26/09/2023 05:12"
$ws.Range("P131").Value = 2.61
$ws.Range("Q131").Value = "23/11/2023 14:13"
$ws.Range("R131").Value = 3.49
$ws.Range("S131").Value = "26/09/2023 05:12"
$ws.Range("T131").Value = 4.02
$ws.Range("U131").Value = "23/11/2023 14:13"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/serbia/prva-liga/metalac-macva-sabac/tIflKsIL/"
